$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 30 / Row 31 swap (Stacks <-> NEARProtocol) ---
$ws.Range("B30").Value = "NEARProtocol"
$ws.Range("C30").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.93"
$ws.Range("E30").Value = "  +12.07%  "

$ws.Range("B31").Value = "Stacks"
$ws.Range("C31").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.88"
$ws.Range("E31").Value = "  +3.46%  "

# --- Price / Volume updates for remaining rows ---
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "66.196.97"
$ws.Range("E2").Value = "  +1.86%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.195.23"
$ws.Range("E3").Value = "  +1.49%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "595.31"
$ws.Range("E5").Value = "  +3.62%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "154.44"
$ws.Range("E6").Value = "  +3.49%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.190.86"
$ws.Range("E8").Value = "  +1.24%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.535"
$ws.Range("E9").Value = "  +1.83%  "
$ws.Range("E10").Value = "  +0.75%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.11"
$ws.Range("E11").Value = "  +0.05%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.514"
$ws.Range("E12").Value = "  +3.30%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000271"
$ws.Range("E13").Value = "  +3.29%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "39.02"
$ws.Range("E14").Value = "  +5.57%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.716.33"
$ws.Range("E15").Value = "  +1.41%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "66.107.21"
$ws.Range("E16").Value = "  +1.63%  "
$ws.Range("E17").Value = "  +4.54%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.198.45"
$ws.Range("E18").Value = "  +0.44%  "
$ws.Range("E19").Value = "  +0.49%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "510.29"
$ws.Range("E20").Value = "  +0.96%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "15.30"
$ws.Range("E21").Value = "  +3.71%  "
$ws.Range("E22").Value = "  +3.45%  "
$ws.Range("E23").Value = "  -0.64%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "7.99"
$ws.Range("E24").Value = "  +3.72%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "84.83"
$ws.Range("E25").Value = "  +0.93%  "
$ws.Range("E26").Value = "  -0.08%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.34"
$ws.Range("E27").Value = "  +5.66%  "
$ws.Range("E28").Value = "  +3.20%  "
$ws.Range("E29").Value = "  +5.45%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "28.33"
$ws.Range("E32").Value = "  +2.87%  "
$ws.Range("E33").Value = "  +3.13%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.00"
$ws.Range("E34").Value = "  +0.28%  "
$ws.Range("E35").Value = "  +0.76%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "54.88"
$ws.Range("E36").Value = "  -0.17%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0902"
$ws.Range("E37").Value = "  +1.07%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "483.95"
$ws.Range("E38").Value = "  +4.18%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0418"
$ws.Range("E39").Value = "  -0.52%  "
$ws.Range("E40").Value = "  -1.90%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "8.85"
$ws.Range("E41").Value = "  +2.48%  "
$ws.Range("E42").Value = "  +4.93%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.299"
$ws.Range("E43").Value = "  +6.38%  "
$ws.Range("E44").Value = "  +12.02%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.931.16"
$ws.Range("E45").Value = "  -4.00%  "
$ws.Range("E46").Value = "  -0.25%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "28.47"
$ws.Range("E47").Value = "  -0.21%  "
$ws.Range("E48").Value = "  -0.02%  "
$ws.Range("E49").Value = "  +2.02%  "
$ws.Range("E50").Value = "  +3.98%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.59"
$ws.Range("E51").Value = "  +6.35%  "
